$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header row
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Inicial"
$ws.Range("C1").Value = "Descripción"

# Clear old row2/row3 leftover data in C column (old numeric votes) just in case
$ws.Range("A2:C9").ClearContents()

# Data rows (party list)
$ws.Range("A3").Value = "Podemos"
$ws.Range("B3").Value = "Podemos"
$ws.Range("C3").Value = "Partido politico"

$ws.Range("A4").Value = "Ciudadanos"
$ws.Range("B4").Value = "Cs"
$ws.Range("C4").Value = "Partido politico"

$ws.Range("A5").Value = "Partido Popular"
$ws.Range("B5").Value = "PP"
$ws.Range("C5").Value = "Partido Politico"

$ws.Range("A6").Value = "Partido Socialista Obrero Español"
$ws.Range("B6").Value = "PSOE"
$ws.Range("C6").Value = "Partido politico"

$ws.Range("A7").Value = "Unión Progreso y Democracia"
$ws.Range("B7").Value = "UpyD"
$ws.Range("C7").Value = "Partido politico"

$ws.Range("A8").Value = "Izquierda Unida"
$ws.Range("B8").Value = "IU"
$ws.Range("C8").Value = "Partido politico"

$ws.Range("A9").Value = "Vox"
$ws.Range("B9").Value = "Vox"
$ws.Range("C9").Value = "Partido politico"

# Column widths (character units, matching the target stored widths of
# 29.5 / 13.1640625 / 16.1640625 as closely as the engine's width
# quantization allows)
$ws.Columns.Item(1).ColumnWidth = 28.6666666666667
$ws.Columns.Item(2).ColumnWidth = 12.3333333333333
$ws.Columns.Item(3).ColumnWidth = 15.3333333333333

# Selection on the active sheet
$ws.Range("A10").Select()
